$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.409254
$ws.Range("H2").Value = 127.227762
$ws.Range("I2").Value = 0.6138221220752584
$ws.Range("J2").Value = 0.6138221220752584
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 173.609943187556
$ws.Range("R2").Value = 1562.489488688004
$ws.Range("S2").Value = 0.09888110750949249
$ws.Range("T2").Value = 0.09888110750949247

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.409254
$ws.Range("H3").Value = 127.227762
$ws.Range("I3").Value = 0.6138221220752584
$ws.Range("J3").Value = 0.6138221220752584
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("Q3").Value = 604.222071941394
$ws.Range("R3").Value = 5437.998647472547
$ws.Range("S3").Value = 0.3441401256073204
$ws.Range("T3").Value = 0.3441401256073204

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.409254
$ws.Range("H4").Value = 127.227762
$ws.Range("I4").Value = 0.6138221220752584
$ws.Range("J4").Value = 0.6138221220752584
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 299.88269119673
$ws.Range("R4").Value = 2698.94422077057
$ws.Range("S4").Value = 0.1708008889584455
$ws.Range("T4").Value = 0.1708008889584455

# Row 5
$ws.Range("I5").Value = 0.07014398987036251
$ws.Range("J5").Value = 0.07014398987036251
$ws.Range("M5").Value = 4.093680666666667
$ws.Range("N5").Value = 12.281042
$ws.Range("O5").Value = 0.1610908176055751
$ws.Range("P5").Value = 0.161090817605575
$ws.Range("Q5").Value = 19.83912547037378
$ws.Range("R5").Value = 178.552129233364
$ws.Range("S5").Value = 0.01129955267833387
$ws.Range("T5").Value = 0.01129955267833387

# Row 6
$ws.Range("I6").Value = 0.07014398987036251
$ws.Range("J6").Value = 0.07014398987036251
$ws.Range("O6").Value = 0.5606512265211691
$ws.Range("P6").Value = 0.5606512265211691
$ws.Range("S6").Value = 0.0393263139539072
$ws.Range("T6").Value = 0.0393263139539072

# Row 7
$ws.Range("I7").Value = 0.07014398987036251
$ws.Range("J7").Value = 0.07014398987036251
$ws.Range("M7").Value = 7.071161666666666
$ws.Range("N7").Value = 21.213485
$ws.Range("O7").Value = 0.2782579558732559
$ws.Range("P7").Value = 0.2782579558732559
$ws.Range("Q7").Value = 34.2688340760411
$ws.Range("R7").Value = 308.4195066843699
$ws.Range("S7").Value = 0.01951812323812144
$ws.Range("T7").Value = 0.01951812323812144

# Row 8
$ws.Range("G8").Value = 21.83492733333334
$ws.Range("H8").Value = 65.50478200000001
$ws.Range("I8").Value = 0.3160338880543792
$ws.Range("J8").Value = 0.3160338880543791
$ws.Range("M8").Value = 4.093680666666667
$ws.Range("N8").Value = 12.281042
$ws.Range("O8").Value = 0.1610908176055751
$ws.Range("P8").Value = 0.161090817605575
$ws.Range("Q8").Value = 89.38521988253824
$ws.Range("R8").Value = 804.466978942844
$ws.Range("S8").Value = 0.05091015741774872
$ws.Range("T8").Value = 0.0509101574177487

# Row 9
$ws.Range("G9").Value = 21.83492733333334
$ws.Range("H9").Value = 65.50478200000001
$ws.Range("I9").Value = 0.3160338880543792
$ws.Range("J9").Value = 0.3160338880543791
$ws.Range("O9").Value = 0.5606512265211691
$ws.Range("P9").Value = 0.5606512265211691
$ws.Range("Q9").Value = 311.091183873134
$ws.Range("R9").Value = 2799.820654858206
$ws.Range("S9").Value = 0.1771847869599415
$ws.Range("T9").Value = 0.1771847869599415

# Row 10
$ws.Range("G10").Value = 21.83492733333334
$ws.Range("H10").Value = 65.50478200000001
$ws.Range("I10").Value = 0.3160338880543792
$ws.Range("J10").Value = 0.3160338880543791
$ws.Range("M10").Value = 7.071161666666666
$ws.Range("N10").Value = 21.213485
$ws.Range("O10").Value = 0.2782579558732559
$ws.Range("P10").Value = 0.2782579558732559
$ws.Range("Q10").Value = 154.3983011539189
$ws.Range("R10").Value = 1389.58471038527
$ws.Range("S10").Value = 0.08793894367668892
$ws.Range("T10").Value = 0.0879389436766889

Write-Output "Updated cells for rows 2-10"